# Apply the "Fixed update to excel issue" change:
#  1. Rename "Requested quantity" header -> "Weekly_PO_Qty" on "Weekly Quantity" sheet
#  2. Rename "Requested quantity" header -> "Monthly_PO_Qty" on "Monthly Trend" sheet
#  3. Add a new "PO Forecast" sheet (after "Monthly Trend") with forecast data

$wb = $excel.ActiveWorkbook
$origActiveSheetName = $wb.ActiveSheet.Name

# --- 1 & 2: header renames -------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: new "PO Forecast" sheet --------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the bold/centered header formatting from an existing header cell so the
# new header reuses the same cell style as the other sheets.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
  @(45340.99999999999,149,21.17114368129885,275.8710087638053),
  @(45389.99999999999,140,5.830337672595747,266.1529536744811),
  @(45403.99999999999,138,4.342076270128183,265.04887152487),
  @(45410.99999999999,136,8.269497358521487,268.0341724657395),
  @(45417.99999999999,135,12.77221792366475,257.6050360055319),
  @(45424.99999999999,134,5.466464381769482,256.8839523663649),
  @(45431.99999999999,132,10.39058455250369,268.6305924947163),
  @(45452.99999999999,128,3.459791526998332,249.4346207973294),
  @(45466.99999999999,126,-0.3757269101645511,247.6711741974744),
  @(45473.99999999999,124,-0.4737774404029909,250.7264529744319),
  @(45480.99999999999,123,-7.528521550317169,261.5550877537919),
  @(45487.99999999999,122,-1.749404978650591,257.8776385182778),
  @(45501.99999999999,119,-3.94257340008563,247.1576619512947),
  @(45508.99999999999,118,-9.799560949866214,252.9143168065282),
  @(45543.99999999999,111,-19.98603266308874,245.4740680115861),
  @(45550.99999999999,110,-13.86626226241932,236.2710259880828),
  @(45557.99999999999,109,-17.85427466029679,239.4138913385822),
  @(45564.99999999999,107,-17.18822484307755,231.0221625681071),
  @(45571.99999999999,106,-24.48159067739511,232.8010878821793),
  @(45578.99999999999,105,-23.27851593201581,225.5293209148228),
  @(45585.99999999999,103,-30.6440741269498,232.5825551840952),
  @(45592.99999999999,102,-29.51248478412736,224.4284844329815),
  @(45599.99999999999,101,-20.6097257917441,231.1621427341103),
  @(45606.99999999999,99,-36.66913334039158,220.0572426043311),
  @(45613.99999999999,98,-32.05551752575744,225.3728245217459),
  @(45620.99999999999,97,-9.959554926799267,230.3457883284937),
  @(45627.99999999999,95,-31.23489272601154,223.3629643322965),
  @(45634.99999999999,94,-29.93220177866176,219.2868836789089),
  @(45641.99999999999,93,-31.78961903673084,220.5449589257795),
  @(45648.99999999999,91,-24.2513142743415,217.6559297773376),
  @(45655.99999999999,90,-40.03838133840313,224.838335799956)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Apply the date/time number format (same style as column A on the other
# sheets) to the "ds" column's data rows.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A32").PasteSpecial(-4122)

# Restore the original active sheet/selection so adding the new sheet doesn't
# change which tab is active.
$wb.Worksheets.Item($origActiveSheetName).Activate()

Write-Host "PO Forecast sheet created and headers renamed."
